$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- Remove the two accounting-entry sheets (no longer needed) ---
$wb.Worksheets("Acc_Disbursement").Delete()
$wb.Worksheets("Acc_Repayment").Delete()

# --- Summary sheet: update the "Over Due" fees row ---
$wsSummary = $wb.Worksheets("Summary")
$wsSummary.Range("B5").Value = 8.88
$wsSummary.Range("E5").Value = 8.88
$wsSummary.Range("F5").Value = 8.88

# --- Repayment schedule: fix the stray placeholder column + correct fee figures ---
$wsSched = $wb.Worksheets("Repayment schedule")

# the empty placeholder cell on the disbursement row was in column P; it belongs in column O
$wsSched.Range("P2").Clear()
$wsSched.Range("N2").Copy()
$wsSched.Range("O2").PasteSpecial(-4122)

$wsSched.Range("J3").Value = 8.8800000000000008
$wsSched.Range("L4").Value = 0
$wsSched.Range("P4").Value = 896.6
$wsSched.Range("J5").Value = 0
$wsSched.Range("K5").Value = 887.72
$wsSched.Range("P5").Value = 887.72

# drop the unused placeholder column values for every repayment installment row
$wsSched.Range("O3:O8").Clear()

# --- Transactions: refresh the entry ids and the fee/balance figures ---
$wsTrans = $wb.Worksheets("Transactions")
$wsTrans.Range("A2").Value = 3189
$wsTrans.Range("I2").Value = 8.8800000000000008
$wsTrans.Range("J2").Value = 4163.24
$wsTrans.Range("A3").Value = 3172
$wsTrans.Range("J6").Clear()

# --- Restore the per-sheet selections & make Transactions the active tab ---
$wb.Worksheets("Input").Range("A2").Select()
$wsSummary.Range("D5").Select()
$wsSched.Range("F7").Select()
$wsTrans.Range("D3").Select()
